# Auto-generated Excel COM-interop script
# Scheduled-runner market refresh: rewrites the live currentAveragePrice* /
# Leve Price / Leve Profit columns (H, I, J, K, L, M, N) for the affected leve
# rows across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59954.92
$ws.Range("I11").Value = 59954.92
$ws.Range("K11").Value = 59954.92
$ws.Range("M11").Value = -59814.92
$ws.Range("H12").Value = 197
$ws.Range("J12").Value = 399.66666
$ws.Range("L12").Value = 399.66666
$ws.Range("N12").Value = -739.66666
$ws.Range("H13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H18").Value = 1200
$ws.Range("I18").Value = 250
$ws.Range("K18").Value = 250
$ws.Range("M18").Value = 34
$ws.Range("H19").Value = 626.9583
$ws.Range("I19").Value = 655.8571
$ws.Range("K19").Value = 655.8571
$ws.Range("M19").Value = -480.8570999999999
$ws.Range("H42").Value = 66.71429
$ws.Range("I42").Value = 66.71429
$ws.Range("K42").Value = 200.14287
$ws.Range("M42").Value = 29.85712999999998
$ws.Range("H70").Value = 3148.3635
$ws.Range("J70").Value = 3013.3
$ws.Range("L70").Value = 9039.900000000001
$ws.Range("N70").Value = -9579.900000000001
$ws.Range("H73").Value = 3148.3635
$ws.Range("J73").Value = 3013.3
$ws.Range("L73").Value = 9039.900000000001
$ws.Range("N73").Value = -10911.9
$ws.Range("H107").Value = 3634.7368
$ws.Range("J107").Value = 6673.625
$ws.Range("L107").Value = 6673.625
$ws.Range("N107").Value = -10513.625
$ws.Range("H127").Value = 5162.8
$ws.Range("I127").Value = 5304
$ws.Range("K127").Value = 15912
$ws.Range("M127").Value = -10952
$ws.Range("H135").Value = 10301.333
$ws.Range("I135").Value = 6928
$ws.Range("K135").Value = 62352
$ws.Range("M135").Value = -59817
$ws.Range("H138").Value = 2952.9487
$ws.Range("J138").Value = 3131.0967
$ws.Range("L138").Value = 9393.2901
$ws.Range("N138").Value = -19673.2901
$ws.Range("H141").Value = 7431.45
$ws.Range("I141").Value = 4187.7856
$ws.Range("K141").Value = 12563.3568
$ws.Range("M141").Value = -7383.356800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 816
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""
$ws.Range("H32").Value = 198340.75
$ws.Range("I32").Value = 258576.5
$ws.Range("K32").Value = 258576.5
$ws.Range("M32").Value = -258289.5
$ws.Range("H63").Value = 5479
$ws.Range("I63").Value = 5131.6665
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 5131.6665
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -4445.6665
$ws.Range("N63").Value = -7372
$ws.Range("H66").Value = 5479
$ws.Range("I66").Value = 5131.6665
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 25658.3325
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -22226.3325
$ws.Range("N66").Value = -36864
$ws.Range("H132").Value = 678071.94
$ws.Range("I132").Value = 1139253.2
$ws.Range("J132").Value = 1672.6666
$ws.Range("K132").Value = 3417759.6
$ws.Range("L132").Value = 5017.9998
$ws.Range("M132").Value = -3415229.6
$ws.Range("N132").Value = -10077.9998
$ws.Range("H134").Value = 49205.332
$ws.Range("J134").Value = 49205.332
$ws.Range("L134").Value = 49205.332
$ws.Range("N134").Value = -59345.332
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1921.5
$ws.Range("I94").Value = 1824.7142
$ws.Range("K94").Value = 1824.7142
$ws.Range("M94").Value = -1373.7142
$ws.Range("H134").Value = 5215782.5
$ws.Range("J134").Value = 33355130
$ws.Range("L134").Value = 100065390
$ws.Range("N134").Value = -100070460
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2142568.8
$ws.Range("I58").Value = 2405.5454
$ws.Range("K58").Value = 2405.5454
$ws.Range("M58").Value = -2202.5454
$ws.Range("H94").Value = 7886.5
$ws.Range("I94").Value = 34582.332
$ws.Range("K94").Value = 34582.332
$ws.Range("M94").Value = -34131.332
$ws.Range("H132").Value = 2594
$ws.Range("J132").Value = 4332.3335
$ws.Range("L132").Value = 12997.0005
$ws.Range("N132").Value = -18057.0005
$ws.Range("H136").Value = 2142568.8
$ws.Range("I136").Value = 2405.5454
$ws.Range("K136").Value = 7216.6362
$ws.Range("M136").Value = -4666.6362
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 5549.5835
$ws.Range("I44").Value = 673.75
$ws.Range("J44").Value = 7987.5
$ws.Range("K44").Value = 2021.25
$ws.Range("L44").Value = 23962.5
$ws.Range("M44").Value = -1623.25
$ws.Range("N44").Value = -24758.5
$ws.Range("H107").Value = 5184.615
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 7222.222
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 21666.666
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -25506.666
$ws.Range("H121").Value = 5559466.5
$ws.Range("I121").Value = 674.25
$ws.Range("J121").Value = 10006500
$ws.Range("K121").Value = 2022.75
$ws.Range("L121").Value = 30019500
$ws.Range("M121").Value = -712.75
$ws.Range("N121").Value = -30022120
$ws.Range("H131").Value = 4574.222
$ws.Range("I131").Value = 1344.5
$ws.Range("K131").Value = 4033.5
$ws.Range("M131").Value = 1006.5
$ws.Range("H134").Value = 5472.3
$ws.Range("I134").Value = 1629.7333
$ws.Range("K134").Value = 4889.199900000001
$ws.Range("M134").Value = 180.8000999999995
$ws.Range("H138").Value = 26066.39
$ws.Range("I138").Value = 41253.375
$ws.Range("J138").Value = 17966.666
$ws.Range("K138").Value = 123760.125
$ws.Range("L138").Value = 53899.99800000001
$ws.Range("M138").Value = -118620.125
$ws.Range("N138").Value = -64179.99800000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 47000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""
$ws.Range("H48").Value = 14000
$ws.Range("I48").Value = 14000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 14000
$ws.Range("M48").Value = -13515
$ws.Range("N48").Value = 0
$ws.Range("H81").Value = 47000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""
$ws.Range("H84").Value = 47000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""
$ws.Range("H100").Value = 49985
$ws.Range("J100").Value = 49985
$ws.Range("L100").Value = 49985
$ws.Range("N100").Value = -52149
$ws.Range("H107").Value = 4195.8667
$ws.Range("I107").Value = 5174.636
$ws.Range("K107").Value = 5174.636
$ws.Range("M107").Value = -3254.636
$ws.Range("H122").Value = 65578.5
$ws.Range("J122").Value = 13171
$ws.Range("L122").Value = 39513
$ws.Range("N122").Value = -44413
$ws.Range("H132").Value = 10966.322
$ws.Range("I132").Value = 8689.046
$ws.Range("J132").Value = 16533
$ws.Range("K132").Value = 26067.138
$ws.Range("L132").Value = 49599
$ws.Range("M132").Value = -23537.138
$ws.Range("N132").Value = -54659
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2923.625
$ws.Range("I40").Value = 1698.3334
$ws.Range("K40").Value = 1698.3334
$ws.Range("M40").Value = -1562.3334
$ws.Range("H132").Value = 9553622
$ws.Range("I132").Value = 22286850
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 66860550
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -66858020
$ws.Range("N132").Value = -16160
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4507.1577
$ws.Range("I81").Value = 4792.8887
$ws.Range("J81").Value = 4250
$ws.Range("K81").Value = 9585.7774
$ws.Range("L81").Value = 8500
$ws.Range("M81").Value = -8524.7774
$ws.Range("N81").Value = -10622
$ws.Range("H84").Value = 4507.1577
$ws.Range("I84").Value = 4792.8887
$ws.Range("J84").Value = 4250
$ws.Range("K84").Value = 47928.887
$ws.Range("L84").Value = 42500
$ws.Range("M84").Value = -42624.887
$ws.Range("N84").Value = -53108
$ws.Range("H132").Value = 8336053
$ws.Range("I132").Value = 8774566
$ws.Range("K132").Value = 26323698
$ws.Range("M132").Value = -26321168
$ws.Range("H136").Value = 5294021.5
$ws.Range("I136").Value = 1061439.1
$ws.Range("K136").Value = 3184317.3
$ws.Range("M136").Value = -3181767.3
